$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which
# stores prices as formatted strings (e.g. thousands-separated, fixed decimals).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.654.70'
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").Value = '3.361.16'
$ws.Range("E3").Value = '  +2.47%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").Value = '531.64'
$ws.Range("E5").Value = '  +2.54%  '

$ws.Range("D6").Value = '173.12'
$ws.Range("E6").Value = '  -4.35%  '

$ws.Range("D7").Value = '0.596'
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").Value = '3.369.36'
$ws.Range("E8").Value = '  +2.83%  '

$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("D10").Value = '0.609'
$ws.Range("E10").Value = '  -1.29%  '

$ws.Range("D11").Value = '53.31'
$ws.Range("E11").Value = '  -7.28%  '

$ws.Range("E12").Value = '  +3.03%  '

$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").Value = '3.874.12'
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").Value = '3.340.07'
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("D18").Value = '17.52'
$ws.Range("E18").Value = '  -0.68%  '

$ws.Range("D19").Value = '63.596.07'
$ws.Range("E19").Value = '  +0.27%  '

$ws.Range("D20").Value = '11.24'
$ws.Range("E20").Value = '  +2.97%  '

$ws.Range("D21").Value = '0.968'
$ws.Range("E21").Value = '  +2.32%  '

$ws.Range("D22").Value = '372.56'
$ws.Range("E22").Value = '  +0.77%  '

$ws.Range("B23").Value = 'RenderToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D23").Value = '11.30'
$ws.Range("E23").Value = '  +0.92%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '81.74'
$ws.Range("E24").Value = '  +2.39%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '4.06'
$ws.Range("E25").Value = '  +6.56%  '

$ws.Range("D26").Value = '3.75'
$ws.Range("E26").Value = '  +1.41%  '

$ws.Range("D27").Value = '6.18'
$ws.Range("E27").Value = '  +3.05%  '

$ws.Range("D28").Value = '2.71'
$ws.Range("E28").Value = '  +3.72%  '

$ws.Range("D29").Value = '11.32'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").Value = '28.91'
$ws.Range("E31").Value = '  +1.77%  '

$ws.Range("D32").Value = '636.48'
$ws.Range("E32").Value = '  -1.83%  '

$ws.Range("D33").Value = '6.45'
$ws.Range("E33").Value = '  -3.01%  '

$ws.Range("D34").Value = '11.21'
$ws.Range("E34").Value = '  +0.90%  '

$ws.Range("E35").Value = '  +1.10%  '

$ws.Range("D36").Value = '57.92'
$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").Value = '36.51'
$ws.Range("E38").Value = '  +1.57%  '

$ws.Range("D39").Value = '0.381'
$ws.Range("E39").Value = '  -1.91%  '

$ws.Range("D40").Value = '0.0₃0726'
$ws.Range("E40").Value = '  +12.13%  '

$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.46%  '

$ws.Range("E42").Value = '  +8.25%  '

$ws.Range("D43").Value = '0.125'
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("D44").Value = '2.937.56'
$ws.Range("E44").Value = '  -1.57%  '

$ws.Range("D45").Value = '3.02'
$ws.Range("E45").Value = '  +8.00%  '

$ws.Range("D46").Value = '2.68'
$ws.Range("E46").Value = '  +3.12%  '

$ws.Range("E47").Value = '  +2.75%  '

$ws.Range("D48").Value = '2.61'
$ws.Range("E48").Value = '  -2.17%  '

$ws.Range("D49").Value = '3.05'
$ws.Range("E49").Value = '  +4.30%  '

$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("D51").Value = '136.26'
$ws.Range("E51").Value = '  +4.68%  '
